$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 1799.875
$ws.Range("I2").Value = 1600
$ws.Range("J2").Value = 1999.75
$ws.Range("K2").Value = 1600
$ws.Range("L2").Value = 1999.75
$ws.Range("M2").Value = -1487
$ws.Range("N2").Value = -2225.75
# row 4
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 56
$ws.Range("K4").Value = 56
$ws.Range("M4").Value = 58
# row 18
$ws.Range("H18").Value = 8142.857
$ws.Range("I18").Value = 1600
$ws.Range("J18").Value = 24500
$ws.Range("K18").Value = 1600
$ws.Range("L18").Value = 24500
$ws.Range("M18").Value = -1316
$ws.Range("N18").Value = -25068
# row 87
$ws.Range("H87").Value = 33323.5
$ws.Range("J87").Value = 33352.91
$ws.Range("L87").Value = 33352.91
$ws.Range("N87").Value = -35848.91
# row 90
$ws.Range("H90").Value = 33323.5
$ws.Range("J90").Value = 33352.91
$ws.Range("L90").Value = 100058.73
$ws.Range("N90").Value = -112538.73
# row 137
$ws.Range("H137").Value = 5527.515
$ws.Range("J137").Value = 7633.75
$ws.Range("L137").Value = 22901.25
$ws.Range("N137").Value = -28001.25
# row 138
$ws.Range("H138").Value = 2524.3691
$ws.Range("J138").Value = 3289
$ws.Range("L138").Value = 9867
$ws.Range("N138").Value = -20147

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 2980872.2
$ws.Range("I32").Value = 477873.9
$ws.Range("K32").Value = 477873.9
$ws.Range("M32").Value = -477586.9
# row 61
$ws.Range("H61").Value = 4355.5
$ws.Range("I61").Value = 4276.926
$ws.Range("J61").Value = 4779.8
$ws.Range("K61").Value = 4276.926
$ws.Range("L61").Value = 4779.8
$ws.Range("M61").Value = -4064.926
$ws.Range("N61").Value = -5203.8
# row 74
$ws.Range("H74").Value = 2221.2334
$ws.Range("I74").Value = 1434.25
$ws.Range("K74").Value = 1434.25
$ws.Range("M74").Value = -560.25
# row 77
$ws.Range("H77").Value = 2221.2334
$ws.Range("I77").Value = 1434.25
$ws.Range("K77").Value = 7171.25
$ws.Range("M77").Value = -2803.25
# row 110
$ws.Range("H110").Value = 1804.5834
$ws.Range("J110").Value = 3133.3333
$ws.Range("L110").Value = 3133.3333
$ws.Range("N110").Value = -7223.3333
# row 132
$ws.Range("H132").Value = 4804.636
$ws.Range("I132").Value = 4928.222
$ws.Range("K132").Value = 14784.666
$ws.Range("M132").Value = -12254.666
# row 136
$ws.Range("H136").Value = 4355.5
$ws.Range("I136").Value = 4276.926
$ws.Range("J136").Value = 4779.8
$ws.Range("K136").Value = 12830.778
$ws.Range("L136").Value = 14339.4
$ws.Range("M136").Value = -10280.778
$ws.Range("N136").Value = -19439.4

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# row 20
$ws.Range("H20").Value = 4870.55
$ws.Range("I20").Value = 4406.727
$ws.Range("K20").Value = 4406.727
$ws.Range("M20").Value = -4159.727
# row 40
$ws.Range("H40").Value = 24000
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
# row 82
$ws.Range("H82").Value = 15274
$ws.Range("I82").Value = 7169
$ws.Range("K82").Value = 7169
$ws.Range("M82").Value = -6786
# row 85
$ws.Range("H85").Value = 15274
$ws.Range("I85").Value = 7169
$ws.Range("K85").Value = 7169
$ws.Range("M85").Value = -5843
# row 86
$ws.Range("H86").Value = 2993.5
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 2993.5
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 2993.5
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -5239.5
# row 89
$ws.Range("H89").Value = 2993.5
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 2993.5
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 14967.5
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -26199.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 2749.1082
$ws.Range("I31").Value = 1629.7
$ws.Range("J31").Value = 4066.0588
$ws.Range("K31").Value = 1629.7
$ws.Range("L31").Value = 4066.0588
$ws.Range("M31").Value = -1334.7
$ws.Range("N31").Value = -4656.0588
# row 34
$ws.Range("H34").Value = 2749.1082
$ws.Range("I34").Value = 1629.7
$ws.Range("J34").Value = 4066.0588
$ws.Range("K34").Value = 1629.7
$ws.Range("L34").Value = 4066.0588
$ws.Range("M34").Value = -1427.7
$ws.Range("N34").Value = -4470.0588
# row 58
$ws.Range("H58").Value = 1000
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
# row 136
$ws.Range("H136").Value = 1000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 4
$ws.Range("H4").Value = 2158889
$ws.Range("I4").Value = 2188680.8
$ws.Range("K4").Value = 6566042.399999999
$ws.Range("M4").Value = -6565930.399999999
# row 5
$ws.Range("H5").Value = 798.1429000000001
$ws.Range("J5").Value = 1150.5294
$ws.Range("L5").Value = 3451.5882
$ws.Range("N5").Value = -3675.5882
# row 92
$ws.Range("H92").Value = 1095.909
$ws.Range("I92").Value = 1302.75
$ws.Range("J92").Value = 977.7143
$ws.Range("K92").Value = 3908.25
$ws.Range("L92").Value = 2933.1429
$ws.Range("M92").Value = -2660.25
$ws.Range("N92").Value = -5429.1429
# row 98
$ws.Range("H98").Value = 186.1
$ws.Range("J98").Value = 177.2
$ws.Range("L98").Value = 531.5999999999999
$ws.Range("N98").Value = -3527.6
# row 135
$ws.Range("H135").Value = 798.1429000000001
$ws.Range("J135").Value = 1150.5294
$ws.Range("L135").Value = 10354.7646
$ws.Range("N135").Value = -15424.7646

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 99
$ws.Range("H99").Value = 7721
$ws.Range("I99").Value = 3164.5
$ws.Range("K99").Value = 3164.5
$ws.Range("M99").Value = -918.5
# row 132
$ws.Range("H132").Value = 9599.071
$ws.Range("I132").Value = 10032.667
$ws.Range("J132").Value = 6997.5
$ws.Range("K132").Value = 30098.001
$ws.Range("L132").Value = 20992.5
$ws.Range("M132").Value = -27568.001
$ws.Range("N132").Value = -26052.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 30030.309
$ws.Range("I7").Value = 36643.445
$ws.Range("K7").Value = 36643.445
$ws.Range("M7").Value = -36531.445
# row 16
$ws.Range("H16").Value = 1475
$ws.Range("I16").Value = 1424.8182
$ws.Range("J16").Value = 1751
$ws.Range("K16").Value = 1424.8182
$ws.Range("L16").Value = 1751
$ws.Range("M16").Value = -1254.8182
$ws.Range("N16").Value = -2091
# row 22
$ws.Range("H22").Value = 1335
$ws.Range("J22").Value = 1548.75
$ws.Range("L22").Value = 1548.75
$ws.Range("N22").Value = -2138.75
# row 27
$ws.Range("H27").Value = 1335
$ws.Range("J27").Value = 1548.75
$ws.Range("L27").Value = 1548.75
$ws.Range("N27").Value = -1762.75
# row 40
$ws.Range("H40").Value = 4030.7778
$ws.Range("I40").Value = 3965.36
$ws.Range("K40").Value = 3965.36
$ws.Range("M40").Value = -3829.36
# row 55
$ws.Range("H55").Value = 325.26086
$ws.Range("I55").Value = 276.5
$ws.Range("J55").Value = 401.1111
$ws.Range("K55").Value = 276.5
$ws.Range("L55").Value = 401.1111
$ws.Range("M55").Value = -103.5
$ws.Range("N55").Value = -747.1111000000001
# row 93
$ws.Range("H93").Value = 38253.11
$ws.Range("I93").Value = 1241.1666
$ws.Range("K93").Value = 1241.1666
$ws.Range("M93").Value = 6.833399999999983
# row 126
$ws.Range("H126").Value = 30030.309
$ws.Range("I126").Value = 36643.445
$ws.Range("K126").Value = 109930.335
$ws.Range("M126").Value = -107460.335
# row 132
$ws.Range("H132").Value = 5997
$ws.Range("I132").Value = 6745
$ws.Range("K132").Value = 20235
$ws.Range("M132").Value = -17705
# row 136
$ws.Range("H136").Value = 3567.32
$ws.Range("I136").Value = 2584.45
$ws.Range("J136").Value = 7498.8
$ws.Range("K136").Value = 7753.349999999999
$ws.Range("L136").Value = 22496.4
$ws.Range("M136").Value = -5203.349999999999
$ws.Range("N136").Value = -27596.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 81
$ws.Range("H81").Value = 10382.286
$ws.Range("J81").Value = 10472
$ws.Range("L81").Value = 20944
$ws.Range("N81").Value = -23066
# row 84
$ws.Range("H84").Value = 10382.286
$ws.Range("J84").Value = 10472
$ws.Range("L84").Value = 104720
$ws.Range("N84").Value = -115328
# row 107
$ws.Range("H107").Value = 50057784
$ws.Range("I107").Value = 960
$ws.Range("J107").Value = 100114610
$ws.Range("K107").Value = 2880
$ws.Range("L107").Value = 300343830
$ws.Range("M107").Value = -960
$ws.Range("N107").Value = -300347670
# row 126
$ws.Range("H126").Value = 3214.7693
$ws.Range("I126").Value = 2754.3333
$ws.Range("K126").Value = 8262.999899999999
$ws.Range("M126").Value = -5792.999899999999
# row 132
$ws.Range("H132").Value = 5812.706
$ws.Range("I132").Value = 5887.8335
$ws.Range("J132").Value = 5249.25
$ws.Range("K132").Value = 17663.5005
$ws.Range("L132").Value = 15747.75
$ws.Range("M132").Value = -15133.5005
$ws.Range("N132").Value = -20807.75
